# User_Interface/test/user_interface_testplan.xlsx
# Commit: "GETs and POSTs for registration and profiles."
#
# The test plan sheet gains two new columns ("Pass Condition" / "Fail
# Condition") inserted right after the "Expected Result" column (D), pushing
# the former Status/Verification/Comments columns two slots to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns (E:F), shifting old E:G -> G:I ------------
$ws.Columns("E:F").Insert()

# Match the width of the neighbouring "Expected Result" column (D).
$ws.Columns("E:F").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# --- Header row -------------------------------------------------------------
$ws.Range("E1").Value = "Pass Condition"
$ws.Range("F1").Value = "Fail Condition"

# --- Data rows: default Pass/Fail condition text ----------------------------
$ws.Range("E2").Value = "All required fields completed"
$ws.Range("F2").Value = "Missing required fields"

$ws.Range("E3").Value = "All required fields completed"
$ws.Range("F3").Value = "Missing required fields"

$ws.Range("E4").Value = "All required fields completed"
$ws.Range("F4").Value = "Missing required fields"

$ws.Range("E5").Value = "All required fields completed"
$ws.Range("F5").Value = "Missing required fields"

$ws.Range("E6").Value = "Deletion success"
$ws.Range("F6").Value = "No module found"

$ws.Range("E7").Value = "All required fields completed"
$ws.Range("F7").Value = "Missing required fields"

$ws.Range("E8").Value = "All required fields completed"
$ws.Range("F8").Value = "Missing required fields"

$ws.Range("E9").Value = "All required fields completed"
$ws.Range("F9").Value = "Missing required fields"
# Row 9's Pass/Fail cells pick up a distinct (black) font colour.
$ws.Range("E9:F9").Font.Color = 0

$ws.Range("E10").Value = "Deletion success"
$ws.Range("F10").Value = "No account found"

# --- View state changes -----------------------------------------------------
$ws.Range("C8").Select()

$wb.Windows.Item(1).WindowState = -4143
